# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计"),
#    populated with the quarter's fund-holding data.
# 2. Update the "总计" (totals) worksheet with a new leading row summarizing
#    the 2022-Q1 data, shifting the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Use "2021-Q1" as a formatting donor: it already carries the bold/bordered
# header style (s=2) and the matching row-index number style used across the
# other quarterly sheets, and has enough rows to cover our 4 data rows.
$srcSheet = $wb.Worksheets.Item("2021-Q1")

# Header row formatting + labels
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row-index column (A) formatting for the 4 data rows
$srcSheet.Range("A2:A5").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)  # xlPasteFormats

# Data rows: code, name, scale, total position, position ratio, market
# value(亿元), rank. Columns B-G are stored as text (matching the source
# data's convention), column H as a number.
$data = @(
  @("001167", "金鹰科技创新股票",     "4.03", "94.55", "5.79", "0.2333", 2),
  @("162102", "金鹰中小盘精选混合",   "4.60", "76.52", "4.82", "0.2217", 2),
  @("210009", "金鹰核心资源混合",     "3.86", "94.96", "5.53", "0.2135", 2),
  @("001613", "长城久祥灵活配置混合", "0.30", "88.54", "4.27", "0.0128", 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $newSheet.Range("A$r").Value = $i

    $newSheet.Range("B$r").Value = "'" + $row[0]
    $newSheet.Range("C$r").Value = "'" + $row[1]
    $newSheet.Range("D$r").Value = "'" + $row[2]
    $newSheet.Range("E$r").Value = "'" + $row[3]
    $newSheet.Range("F$r").Value = "'" + $row[4]
    $newSheet.Range("G$r").Value = "'" + $row[5]
    # Leading apostrophe forces text entry (avoids "4.03" -> number
    # coercion); reset the style afterwards so the quote-prefix flag
    # doesn't leave a stray style on the cell.
    $newSheet.Range("B$r`:G$r").Style = "Normal"

    $newSheet.Range("H$r").Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: update "总计" with a new first data row for 2022-Q1
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Re-apply the row-index number style (s=2) to the new A2, then clear any
# inherited formatting from B2:D2 so they match the plain data cells below.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.68
$totalSheet.Range("B2:D2").Style = "Normal"

# The A column holds a literal 0-based row counter (not a formula), so the
# rows pushed down by the insert need to be renumbered by hand.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
